$d = $word.ActiveDocument

# --- Step 1: append two trailing spaces to the first run (keeps default/no-color formatting) ---
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$pos1 = $r.End - 1
$ip1 = $d.Range($pos1, $pos1)
$ip1.InsertAfter("  ")

# --- Step 2: append the red "(This is a change – Version for branch " run ---
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$pos2 = $r.End - 1
$text2 = "(This is a change – Version for branch "
$ip2 = $d.Range($pos2, $pos2)
$ip2.InsertAfter($text2)
$run2 = $d.Range($pos2, $pos2 + $text2.Length)
$run2.Font.Color = 192

# --- Step 3: append the red "main" run ---
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$pos3 = $r.End - 1
$text3 = "main"
$ip3 = $d.Range($pos3, $pos3)
$ip3.InsertAfter($text3)
$run3 = $d.Range($pos3, $pos3 + $text3.Length)
$run3.Font.Color = 192

# --- Step 4: append the red ")" run ---
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$pos4 = $r.End - 1
$text4 = ")"
$ip4 = $d.Range($pos4, $pos4)
$ip4.InsertAfter($text4)
$run4 = $d.Range($pos4, $pos4 + $text4.Length)
$run4.Font.Color = 192

# --- Step 5: insert a brand-new, fully empty paragraph right after paragraph 1 ---
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$ipPara = $d.Range($r.End, $r.End)
$null = $ipPara.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
